# Applies the "Ref_clean electricity" BC edit described by the commit:
#   - re-orders the Technology values (col F) for the Base/Shoulder load
#     "Market share_class" blocks (rows 4-16 and 18-23)
#   - bumps the Market share_class_min year values in row 3 / O17:W17 /
#     O24:W24 from 1 -> 0.999
#   - widens columns A, E, F, G to fit their (new, longer) contents
#   - updates the sheet's selection/scroll state

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Activate()

# ---------------------------------------------------------------------
# 1. Row 3 / 17 / 24 "Market share_class_min" series: 1 -> 0.999
# ---------------------------------------------------------------------
$ws.Range("N3:W3").Value = 0.999
$ws.Range("O17:W17").Value = 0.999
$ws.Range("O24:W24").Value = 0.999

# ---------------------------------------------------------------------
# 2. Column F (Technology) re-ordering
# ---------------------------------------------------------------------
# Base Load block (rows 4-16): Hydro River, Solar PV, Wind Onshore,
# Wind Offshore, Nuclear, Nuclear SMR, Hydro Reservoir, Solar PV + Battery,
# Wind + Battery, Wind + Seasonal, Biomass, Biomass CCS, Geoexchange
# becomes:
$ws.Range("F4").Value2  = "Biomass"
$ws.Range("F5").Value2  = "Biomass CCS"
$ws.Range("F6").Value2  = "Geoexchange"
$ws.Range("F7").Value2  = "Nuclear"
$ws.Range("F8").Value2  = "Nuclear SMR"
$ws.Range("F9").Value2  = "Hydro Reservoir"
$ws.Range("F10").Value2 = "Hydro River"
$ws.Range("F11").Value2 = "Solar PV"
$ws.Range("F12").Value2 = "Solar PV + Battery"
$ws.Range("F13").Value2 = "Wind Onshore"
$ws.Range("F14").Value2 = "Wind + Battery"
$ws.Range("F15").Value2 = "Wind + Seasonal"
$ws.Range("F16").Value2 = "Wind Offshore"

# Shoulder Load block (rows 18-23): Hydro Reservoir, Solar PV + Battery,
# Wind + Battery, Wind + Seasonal, Biomass, Biomass CCS
# becomes:
$ws.Range("F18").Value2 = "Biomass"
$ws.Range("F19").Value2 = "Biomass CCS"
$ws.Range("F20").Value2 = "Hydro Reservoir"
$ws.Range("F21").Value2 = "Solar PV + Battery"
$ws.Range("F22").Value2 = "Wind + Battery"
$ws.Range("F23").Value2 = "Wind + Seasonal"

# ---------------------------------------------------------------------
# 3. Column widths for A, E, F, G (best-fit to new content)
# ---------------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 50.833333333333336
$ws.Columns.Item(5).ColumnWidth = 12.833333333333334
$ws.Columns.Item(6).ColumnWidth = 15.5
$ws.Columns.Item(7).ColumnWidth = 22.166666666666668

# ---------------------------------------------------------------------
# 4. Selection / view state
# ---------------------------------------------------------------------
$ws.Range("A3:W26").Select() | Out-Null
$excel.ActiveWindow.ScrollColumn = 2
